# Week 15 simulation data update.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: per-play yardage logs get Week 15 plays appended.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$cur = $ws.Range("B2").Text
$ws.Range("B2").Value = $cur + " 4 3 0 4 0 5 5 4 0 2 2 7 4 5 4 11 1 4 10 10 6 1 3 2 3 3 4 3 4 1 5 5"

$cur = $ws.Range("B3").Text
$ws.Range("B3").Value = $cur + " 4 18 7 15 3 12 24 7 4 0 12 13 9 14 20 12 5 11 5 10 2 4"

$cur = $ws.Range("C2").Text
$ws.Range("C2").Value = $cur + " 2 0 1 3 3 11 2 4 6 5 0 13 6 5 12 8 4 6 7 2 0"

$cur = $ws.Range("C3").Text
$ws.Range("C3").Value = $cur + " 14 5 12 2 3 19 9 43 -1 6 10 8 19 19 7"

# ---------------------------------------------------------------------
# OFF sheet: season totals through Week 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 178
$ws.Range("F2").Value = 48
$ws.Range("G2").Value = 43
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 24
$ws.Range("N2").Value = 16

$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 142
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 108
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 48
$ws.Range("J3").Value = 51
$ws.Range("L3").Value = 301
$ws.Range("M3").Value = 204
$ws.Range("Q3").Value = 498

# ---------------------------------------------------------------------
# DEF sheet: season totals through Week 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 148
$ws.Range("F2").Value = 46
$ws.Range("G2").Value = 36
$ws.Range("J2").Value = 19
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 15
$ws.Range("P2").Value = 9

$ws.Range("C3").Value = 157
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 86
$ws.Range("G3").Value = 29
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 47
$ws.Range("J3").Value = 44
$ws.Range("L3").Value = 259
$ws.Range("M3").Value = 158
$ws.Range("Q3").Value = 435

# ---------------------------------------------------------------------
# ST sheet: season totals + per-game logs through Week 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 79
$ws.Range("D2").Value = 50
$ws.Range("F2").Value = 340
$ws.Range("G2").Value = 328
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 168
$ws.Range("K2").Value = 159

$ws.Range("B3").Value = 66

$cur = $ws.Range("B4").Text
$ws.Range("B4").Value = $cur + " 66 66"

$cur = $ws.Range("D3").Text
$ws.Range("D3").Value = $cur + " 43 31 59 35 50 50"

$cur = $ws.Range("B5").Text
$ws.Range("B5").Value = $cur + " 28 20"

$cur = $ws.Range("D4").Text
$ws.Range("D4").Value = $cur + " 0 0 0 0 0 0"

$cur = $ws.Range("D5").Text
$ws.Range("D5").Value = $cur + " 0 4 0 1 0 0 0"

$cur = $ws.Range("B6").Text
$ws.Range("B6").Value = $cur + " 21 15"

# ---------------------------------------------------------------------
# TURNS sheet: season totals through Week 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: season totals through Week 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("D2").Value = 9
$ws.Range("B3").Value = 25
